$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update PID tuning input values (Left Motor = column B, Right Motor = column D) ---

# T (row 30)
$ws.Range("B30").Value = 0.05
$ws.Range("D30").Value = 0.05

# TS (row 31)
$ws.Range("B31").Value = 0.12
$ws.Range("D31").Value = 0.1

# y1 (row 34)
$ws.Range("B34").Value = 70.91
$ws.Range("D34").Value = 77.67

# y2 (row 35)
$ws.Range("B35").Value = 87.98
$ws.Range("D35").Value = 95.28

# --- Update the view state: scroll/selection moved to D36, zoom bumped to 164% ---
[void]$ws.Range("D36").Select()
$excel.ActiveWindow.Zoom = 164
